$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D range to text format before writing, so numeric-looking
# strings like "1.147.72" are preserved verbatim as text (not parsed as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.589.31"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "1.794.66"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "231.64"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").Value = "0.5906"
$ws.Range("E6").Value = "  -1.05%  "

$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "0.2778"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("D9").Value = "23.44"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").Value = "0.06768"
$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("D11").Value = "0.07558"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").Value = "1.796.05"
$ws.Range("E12").Value = "  -2.15%  "

$ws.Range("D13").Value = "4.802"
$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").Value = "0.6144"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "2.037.22"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "75.80"
$ws.Range("E16").Value = "  -3.36%  "

$ws.Range("D17").Value = "0.000008932"
$ws.Range("E17").Value = "  -8.11%  "

$ws.Range("D18").Value = "28.569.72"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("D19").Value = "5.430"
$ws.Range("E19").Value = "  -5.34%  "

$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "209.29"
$ws.Range("E21").Value = "  -5.65%  "

$ws.Range("D22").Value = "11.49"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").Value = "6.841"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").Value = "152.43"
$ws.Range("E25").Value = "  -2.26%  "

$ws.Range("D26").Value = "8.034"
$ws.Range("E26").Value = "  +1.10%  "

$ws.Range("D27").Value = "0.1263"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").Value = "16.45"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "1.415"
$ws.Range("E29").Value = "  -2.28%  "

$ws.Range("D30").Value = "0.06168"
$ws.Range("E30").Value = "  -7.43%  "

$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "3.794"
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("D35").Value = "1.051"
$ws.Range("E35").Value = "  -3.49%  "

$ws.Range("D36").Value = "0.6422"
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("D37").Value = "2.501"
$ws.Range("E37").Value = "  -1.77%  "

$ws.Range("D38").Value = "2.712"
$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("D39").Value = "0.01692"
$ws.Range("E39").Value = "  -2.32%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.147.72"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.335"
$ws.Range("E41").Value = "  -2.49%  "

$ws.Range("D42").Value = "0.8737"
$ws.Range("E42").Value = "  -3.05%  "

$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "100.48"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "1.944.88"
$ws.Range("E45").Value = "  -1.86%  "

$ws.Range("D46").Value = "60.25"
$ws.Range("E46").Value = "  -2.57%  "

$ws.Range("E47").Value = "  -3.62%  "

$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").Value = "8.369"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").Value = "0.05456"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("D51").Value = "0.4473"
$ws.Range("E51").Value = "  -1.87%  "

# Restore default style on column D (clears the temporary text number format)
$ws.Range("D2:D51").Style = "Normal"
